# Trade #43 (internal trade id 71, EMAArbitrage strategy) closes out, and a
# new trade (internal trade id 100, MarketMaking strategy) is opened.
# This updates the aggregate Summary / Strategy Status sheets, the closed
# trade's row on "All Trades" + "EMAArbitrage", and appends the freshly
# opened trade's row to "All Trades" + "MarketMaking".

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $addr, $text) {
    # Force text storage so date/time-shaped strings ("2026-02-18",
    # "00:17:26") aren't auto-converted into date/time serial numbers.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.42   # Current Capital
$summary.Range("B4").Value = 0.52      # Total P&L $
$summary.Range("B5").Value = 0.15      # Total P&L %
$summary.Range("B6").Value = 71        # Total Trades
$summary.Range("B8").Value = 29        # Losing Trades
$summary.Range("B9").Value = 50.7      # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - EMAArbitrage row (row 2)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C2").Value = 100.27
$status.Range("D2").Value = 3
$status.Range("E2").Value = 0.27
$status.Range("F2").Value = 0.27
$status.Range("G2").Value = 66.67

# ---------------------------------------------------------------------
# All Trades sheet - close out trade #71 (row 72) and append new trade
# #100 (row 101)
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Range("G72").Value = 0.31
$allTrades.Range("H72").Value = "CLOSED"
$allTrades.Range("I72").Value = -13.8889
$allTrades.Range("J72").Value = -0.05
$allTrades.Range("K72").Value = 100.27
$allTrades.Range("L72").Value = "early_exit"
$allTrades.Range("M72").Value = 0.1

$allTrades.Range("A101").Value = 100
Set-TextCell $allTrades "B101" "2026-02-18"
Set-TextCell $allTrades "C101" "00:17:26"
$allTrades.Range("D101").Value = "MarketMaking"
$allTrades.Range("E101").Value = "DOWN"
$allTrades.Range("F101").Value = 0.36
$allTrades.Range("G101").Value = ""
$allTrades.Range("H101").Value = "OPEN"
$allTrades.Range("I101").Value = 0
$allTrades.Range("J101").Value = 0
$allTrades.Range("K101").Value = 99.410254715139
$allTrades.Range("L101").Value = ""
$allTrades.Range("M101").Value = 0
$allTrades.Range("N101").Value = 0
$allTrades.Range("O101").Value = 0
$allTrades.Range("P101").Value = 0.6
$allTrades.Range("Q101").Value = "Normal spread capture: 198 bps"

# ---------------------------------------------------------------------
# EMAArbitrage sheet - close out trade #71 (row 4)
# ---------------------------------------------------------------------
$emaArb = $wb.Worksheets.Item("EMAArbitrage")
$emaArb.Range("G4").Value = 0.31
$emaArb.Range("H4").Value = "CLOSED"
$emaArb.Range("I4").Value = -13.8889
$emaArb.Range("J4").Value = -0.05
$emaArb.Range("K4").Value = 100.27
$emaArb.Range("P4").Value = "early_exit"
$emaArb.Range("Q4").Value = 0.1

# ---------------------------------------------------------------------
# MarketMaking sheet - append new trade #100 (row 33)
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Range("A33").Value = 100
Set-TextCell $marketMaking "B33" "2026-02-18"
Set-TextCell $marketMaking "C33" "00:17:26"
$marketMaking.Range("D33").Value = "MarketMaking"
$marketMaking.Range("E33").Value = "DOWN"
$marketMaking.Range("F33").Value = 0.36
$marketMaking.Range("G33").Value = ""
$marketMaking.Range("H33").Value = "OPEN"
$marketMaking.Range("I33").Value = 0
$marketMaking.Range("J33").Value = 0
$marketMaking.Range("K33").Value = 99.410254715139
$marketMaking.Range("L33").Value = 0
$marketMaking.Range("M33").Value = 0
$marketMaking.Range("N33").Value = 0.6
$marketMaking.Range("O33").Value = "Normal spread capture: 198 bps"
$marketMaking.Range("P33").Value = ""
$marketMaking.Range("Q33").Value = 0
